# Update countries & provincias Spain
# Refresh COVID case statistics table and fix alphabetical/value
# ordering of a few tied-rank countries (Botsuana/Malaui/Suazilandia,
# Montserrat/Islas Malvinas), plus bump the "last updated" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Range("A1").Value = "Datos actualizados a 23 de Octubre de 2020 a las 01:49"

# Row 4
$ws.Range("B4").Value = 8654232
$ws.Range("C4").Value = 68882
$ws.Range("D4").Value = 5647254
$ws.Range("E4").Value = 2778689
$ws.Range("G4").Value = 881
$ws.Range("H4").Value = 228289

# Row 6
$ws.Range("D6").Value = 4785297
$ws.Range("E6").Value = 391375
$ws.Range("G6").Value = 503
$ws.Range("H6").Value = 155962

# Row 9
$ws.Range("B9").Value = 1053650
$ws.Range("C9").Value = 16325
$ws.Range("D9").Value = 851854
$ws.Range("E9").Value = 173839
$ws.Range("G9").Value = 438
$ws.Range("H9").Value = 27957

# Row 31
$ws.Range("B31").Value = 223065
$ws.Range("C31").Value = 14150
$ws.Range("D31").Value = 87225
$ws.Range("E31").Value = 133995
$ws.Range("G31").Value = 106
$ws.Range("H31").Value = 1845

# Row 40
$ws.Range("B40").Value = 127227
$ws.Range("C40").Value = 792
$ws.Range("D40").Value = 103398
$ws.Range("E40").Value = 21217
$ws.Range("G40").Value = 15
$ws.Range("H40").Value = 2612

# Row 51
$ws.Range("E51").Value = 39273
$ws.Range("G51").Value = 7
$ws.Range("H51").Value = 2046

# Row 58
$ws.Range("B58").Value = 79211
$ws.Range("C58").Value = 304
$ws.Range("D58").Value = 75840
$ws.Range("E58").Value = 3063

# Row 125
$ws.Range("A125").Value = "Botsuana"
$ws.Range("B125").Value = 5923
$ws.Range("C125").Value = 314
$ws.Range("D125").Value = 927
$ws.Range("E125").Value = 4975
$ws.Range("H125").Value = 21

# Row 126
$ws.Range("A126").Value = "Malaui"
$ws.Range("B126").Value = 5874
$ws.Range("C126").Value = 10
$ws.Range("D126").Value = 4764
$ws.Range("E126").Value = 927
$ws.Range("H126").Value = 183

# Row 127
$ws.Range("A127").Value = "Suazilandia"
$ws.Range("B127").Value = 5814
$ws.Range("C127").Value = 9
$ws.Range("D127").Value = 5468
$ws.Range("E127").Value = 230
$ws.Range("H127").Value = 116

# Row 132
$ws.Range("B132").Value = 5267
$ws.Range("C132").Value = 43
$ws.Range("D132").Value = 1655
$ws.Range("E132").Value = 3352
$ws.Range("G132").Value = 3
$ws.Range("H132").Value = 260

# Row 135
$ws.Range("B135").Value = 5154
$ws.Range("C135").Value = 4
$ws.Range("D135").Value = 4995

# Row 137
$ws.Range("D137").Value = 4961
$ws.Range("E137").Value = 30

# Row 140
$ws.Range("B140").Value = 4862
$ws.Range("C140").Value = 4
$ws.Range("E140").Value = 2876

# Row 147
$ws.Range("B147").Value = 3877
$ws.Range("C147").Value = 27
$ws.Range("D147").Value = 2853
$ws.Range("E147").Value = 907
$ws.Range("G147").Value = 1
$ws.Range("H147").Value = 117

# Row 155
$ws.Range("B155").Value = 2701
$ws.Range("C155").Value = 38
$ws.Range("D155").Value = 2204
$ws.Range("E155").Value = 444

# Row 216
$ws.Range("A216").Value = "Montserrat"
$ws.Range("D216").Value = 12
$ws.Range("H216").Value = 1

# Row 217
$ws.Range("A217").Value = "Islas Malvinas"
$ws.Range("D217").Value = 13
$ws.Range("H217").Value = 0

